$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set cell values while cells are still at default style (so quote-prefixed
# numeric-looking text strings are stored as literal Text, not parsed as dates/numbers).

$ws.Cells.Item(162, 1).Value = "Year 5"
$ws.Cells.Item(162, 2).Value = "B2-C1"
$ws.Cells.Item(162, 3).Value = "endocrinology"
$ws.Cells.Item(162, 4).Value = "'1"
$ws.Cells.Item(162, 5).Value = "'03/01/2026"
$ws.Cells.Item(162, 6).Value = "'09:00:00"
$ws.Cells.Item(162, 7).Value = 360

$ws.Cells.Item(163, 1).Value = "Year 5"
$ws.Cells.Item(163, 2).Value = "B2-C1"
$ws.Cells.Item(163, 3).Value = "endocrinology"
$ws.Cells.Item(163, 4).Value = "'2"
$ws.Cells.Item(163, 5).Value = "'04/01/2026"
$ws.Cells.Item(163, 6).Value = "'09:00:00"
$ws.Cells.Item(163, 7).Value = 360

$ws.Cells.Item(164, 1).Value = "Year 5"
$ws.Cells.Item(164, 2).Value = "B2-C1"
$ws.Cells.Item(164, 3).Value = "endocrinology"
$ws.Cells.Item(164, 4).Value = "'3"
$ws.Cells.Item(164, 5).Value = "'05/01/2026"
$ws.Cells.Item(164, 6).Value = "'09:00:00"
$ws.Cells.Item(164, 7).Value = 360

$ws.Cells.Item(165, 1).Value = "Year 5"
$ws.Cells.Item(165, 2).Value = "B2-C1"
$ws.Cells.Item(165, 3).Value = "endocrinology"
$ws.Cells.Item(165, 4).Value = "'4"
$ws.Cells.Item(165, 5).Value = "'06/01/2026"
$ws.Cells.Item(165, 6).Value = "'09:00:00"
$ws.Cells.Item(165, 7).Value = 360

$ws.Cells.Item(166, 1).Value = "Year 5"
$ws.Cells.Item(166, 2).Value = "B2-C1"
$ws.Cells.Item(166, 3).Value = "endocrinology"
$ws.Cells.Item(166, 4).Value = "'5"
$ws.Cells.Item(166, 5).Value = "'07/01/2026"
$ws.Cells.Item(166, 6).Value = "'09:00:00"
$ws.Cells.Item(166, 7).Value = 360

$ws.Cells.Item(167, 1).Value = "Year 5"
$ws.Cells.Item(167, 2).Value = "B2-C1"
$ws.Cells.Item(167, 3).Value = "endocrinology"
$ws.Cells.Item(167, 4).Value = "'6"
$ws.Cells.Item(167, 5).Value = "'10/01/2026"
$ws.Cells.Item(167, 6).Value = "'09:00:00"
$ws.Cells.Item(167, 7).Value = 360

$ws.Cells.Item(168, 1).Value = "Year 5"
$ws.Cells.Item(168, 2).Value = "B2-C1"
$ws.Cells.Item(168, 3).Value = "endocrinology"
$ws.Cells.Item(168, 4).Value = "'7"
$ws.Cells.Item(168, 5).Value = "'11/01/2026"
$ws.Cells.Item(168, 6).Value = "'09:00:00"
$ws.Cells.Item(168, 7).Value = 360

$ws.Cells.Item(169, 1).Value = "Year 5"
$ws.Cells.Item(169, 2).Value = "B2-C1"
$ws.Cells.Item(169, 3).Value = "endocrinology"
$ws.Cells.Item(169, 4).Value = "'8"
$ws.Cells.Item(169, 5).Value = "'12/01/2026"
$ws.Cells.Item(169, 6).Value = "'09:00:00"
$ws.Cells.Item(169, 7).Value = 360

$ws.Cells.Item(170, 1).Value = "Year 5"
$ws.Cells.Item(170, 2).Value = "B2-C1"
$ws.Cells.Item(170, 3).Value = "endocrinology"
$ws.Cells.Item(170, 4).Value = "'9"
$ws.Cells.Item(170, 5).Value = "'13/01/2026"
$ws.Cells.Item(170, 6).Value = "'09:00:00"
$ws.Cells.Item(170, 7).Value = 360

$ws.Cells.Item(171, 1).Value = "Year 5"
$ws.Cells.Item(171, 2).Value = "B2-C1"
$ws.Cells.Item(171, 3).Value = "endocrinology"
$ws.Cells.Item(171, 4).Value = "'10"
$ws.Cells.Item(171, 5).Value = "'14/01/2026"
$ws.Cells.Item(171, 6).Value = "'09:00:00"
$ws.Cells.Item(171, 7).Value = 360

$ws.Cells.Item(172, 1).Value = "Year 5"
$ws.Cells.Item(172, 2).Value = "B2-C1"
$ws.Cells.Item(172, 3).Value = "gastroenterology"
$ws.Cells.Item(172, 4).Value = "'1"
$ws.Cells.Item(172, 5).Value = "'17/01/2026"
$ws.Cells.Item(172, 6).Value = "'09:00:00"
$ws.Cells.Item(172, 7).Value = 360

$ws.Cells.Item(173, 1).Value = "Year 5"
$ws.Cells.Item(173, 2).Value = "B2-C1"
$ws.Cells.Item(173, 3).Value = "gastroenterology"
$ws.Cells.Item(173, 4).Value = "'2"
$ws.Cells.Item(173, 5).Value = "'18/01/2026"
$ws.Cells.Item(173, 6).Value = "'09:00:00"
$ws.Cells.Item(173, 7).Value = 360

$ws.Cells.Item(174, 1).Value = "Year 5"
$ws.Cells.Item(174, 2).Value = "B2-C1"
$ws.Cells.Item(174, 3).Value = "gastroenterology"
$ws.Cells.Item(174, 4).Value = "'3"
$ws.Cells.Item(174, 5).Value = "'19/01/2026"
$ws.Cells.Item(174, 6).Value = "'09:00:00"
$ws.Cells.Item(174, 7).Value = 360

$ws.Cells.Item(175, 1).Value = "Year 5"
$ws.Cells.Item(175, 2).Value = "B2-C1"
$ws.Cells.Item(175, 3).Value = "gastroenterology"
$ws.Cells.Item(175, 4).Value = "'4"
$ws.Cells.Item(175, 5).Value = "'20/01/2026"
$ws.Cells.Item(175, 6).Value = "'09:00:00"
$ws.Cells.Item(175, 7).Value = 360

$ws.Cells.Item(176, 1).Value = "Year 5"
$ws.Cells.Item(176, 2).Value = "B2-C1"
$ws.Cells.Item(176, 3).Value = "gastroenterology"
$ws.Cells.Item(176, 4).Value = "'5"
$ws.Cells.Item(176, 5).Value = "'21/01/2026"
$ws.Cells.Item(176, 6).Value = "'09:00:00"
$ws.Cells.Item(176, 7).Value = 360

$ws.Cells.Item(177, 1).Value = "Year 5"
$ws.Cells.Item(177, 2).Value = "B2-C1"
$ws.Cells.Item(177, 3).Value = "gastroenterology"
$ws.Cells.Item(177, 4).Value = "'6"
$ws.Cells.Item(177, 5).Value = "'07/02/2026"
$ws.Cells.Item(177, 6).Value = "'09:00:00"
$ws.Cells.Item(177, 7).Value = 360

$ws.Cells.Item(178, 1).Value = "Year 5"
$ws.Cells.Item(178, 2).Value = "B2-C1"
$ws.Cells.Item(178, 3).Value = "gastroenterology"
$ws.Cells.Item(178, 4).Value = "'7"
$ws.Cells.Item(178, 5).Value = "'08/02/2026"
$ws.Cells.Item(178, 6).Value = "'09:00:00"
$ws.Cells.Item(178, 7).Value = 360

$ws.Cells.Item(179, 1).Value = "Year 5"
$ws.Cells.Item(179, 2).Value = "B2-C1"
$ws.Cells.Item(179, 3).Value = "gastroenterology"
$ws.Cells.Item(179, 4).Value = "'8"
$ws.Cells.Item(179, 5).Value = "'09/02/2026"
$ws.Cells.Item(179, 6).Value = "'09:00:00"
$ws.Cells.Item(179, 7).Value = 360

$ws.Cells.Item(180, 1).Value = "Year 5"
$ws.Cells.Item(180, 2).Value = "B2-C1"
$ws.Cells.Item(180, 3).Value = "gastroenterology"
$ws.Cells.Item(180, 4).Value = "'9"
$ws.Cells.Item(180, 5).Value = "'10/02/2026"
$ws.Cells.Item(180, 6).Value = "'09:00:00"
$ws.Cells.Item(180, 7).Value = 360

$ws.Cells.Item(181, 1).Value = "Year 5"
$ws.Cells.Item(181, 2).Value = "B2-C1"
$ws.Cells.Item(181, 3).Value = "gastroenterology"
$ws.Cells.Item(181, 4).Value = "'10"
$ws.Cells.Item(181, 5).Value = "'11/02/2026"
$ws.Cells.Item(181, 6).Value = "'09:00:00"
$ws.Cells.Item(181, 7).Value = 360

$ws.Cells.Item(182, 1).Value = "Year 5"
$ws.Cells.Item(182, 2).Value = "B2-C1"
$ws.Cells.Item(182, 3).Value = "nephrology"
$ws.Cells.Item(182, 4).Value = "'1"
$ws.Cells.Item(182, 5).Value = "'13/12/2025"
$ws.Cells.Item(182, 6).Value = "'09:00:00"
$ws.Cells.Item(182, 7).Value = 360

$ws.Cells.Item(183, 1).Value = "Year 5"
$ws.Cells.Item(183, 2).Value = "B2-C1"
$ws.Cells.Item(183, 3).Value = "nephrology"
$ws.Cells.Item(183, 4).Value = "'2"
$ws.Cells.Item(183, 5).Value = "'14/12/2025"
$ws.Cells.Item(183, 6).Value = "'09:00:00"
$ws.Cells.Item(183, 7).Value = 360

$ws.Cells.Item(184, 1).Value = "Year 5"
$ws.Cells.Item(184, 2).Value = "B2-C1"
$ws.Cells.Item(184, 3).Value = "nephrology"
$ws.Cells.Item(184, 4).Value = "'3"
$ws.Cells.Item(184, 5).Value = "'15/12/2025"
$ws.Cells.Item(184, 6).Value = "'09:00:00"
$ws.Cells.Item(184, 7).Value = 360

$ws.Cells.Item(185, 1).Value = "Year 5"
$ws.Cells.Item(185, 2).Value = "B2-C1"
$ws.Cells.Item(185, 3).Value = "nephrology"
$ws.Cells.Item(185, 4).Value = "'4"
$ws.Cells.Item(185, 5).Value = "'16/12/2025"
$ws.Cells.Item(185, 6).Value = "'09:00:00"
$ws.Cells.Item(185, 7).Value = 360

$ws.Cells.Item(186, 1).Value = "Year 5"
$ws.Cells.Item(186, 2).Value = "B2-C1"
$ws.Cells.Item(186, 3).Value = "nephrology"
$ws.Cells.Item(186, 4).Value = "'5"
$ws.Cells.Item(186, 5).Value = "'17/12/2025"
$ws.Cells.Item(186, 6).Value = "'09:00:00"
$ws.Cells.Item(186, 7).Value = 360

$ws.Cells.Item(187, 1).Value = "Year 5"
$ws.Cells.Item(187, 2).Value = "B2-C1"
$ws.Cells.Item(187, 3).Value = "neurology"
$ws.Cells.Item(187, 4).Value = "'1"
$ws.Cells.Item(187, 5).Value = "'20/12/2025"
$ws.Cells.Item(187, 6).Value = "'09:00:00"
$ws.Cells.Item(187, 7).Value = 360

$ws.Cells.Item(188, 1).Value = "Year 5"
$ws.Cells.Item(188, 2).Value = "B2-C1"
$ws.Cells.Item(188, 3).Value = "neurology"
$ws.Cells.Item(188, 4).Value = "'2"
$ws.Cells.Item(188, 5).Value = "'21/12/2025"
$ws.Cells.Item(188, 6).Value = "'09:00:00"
$ws.Cells.Item(188, 7).Value = 360

$ws.Cells.Item(189, 1).Value = "Year 5"
$ws.Cells.Item(189, 2).Value = "B2-C1"
$ws.Cells.Item(189, 3).Value = "neurology"
$ws.Cells.Item(189, 4).Value = "'3"
$ws.Cells.Item(189, 5).Value = "'22/12/2025"
$ws.Cells.Item(189, 6).Value = "'09:00:00"
$ws.Cells.Item(189, 7).Value = 360

$ws.Cells.Item(190, 1).Value = "Year 5"
$ws.Cells.Item(190, 2).Value = "B2-C1"
$ws.Cells.Item(190, 3).Value = "neurology"
$ws.Cells.Item(190, 4).Value = "'4"
$ws.Cells.Item(190, 5).Value = "'23/12/2025"
$ws.Cells.Item(190, 6).Value = "'09:00:00"
$ws.Cells.Item(190, 7).Value = 360

$ws.Cells.Item(191, 1).Value = "Year 5"
$ws.Cells.Item(191, 2).Value = "B2-C1"
$ws.Cells.Item(191, 3).Value = "neurology"
$ws.Cells.Item(191, 4).Value = "'5"
$ws.Cells.Item(191, 5).Value = "'27/12/2025"
$ws.Cells.Item(191, 6).Value = "'09:00:00"
$ws.Cells.Item(191, 7).Value = 360

$ws.Cells.Item(192, 1).Value = "Year 5"
$ws.Cells.Item(192, 2).Value = "B2-C1"
$ws.Cells.Item(192, 3).Value = "neurology"
$ws.Cells.Item(192, 4).Value = "'6"
$ws.Cells.Item(192, 5).Value = "'28/12/2025"
$ws.Cells.Item(192, 6).Value = "'09:00:00"
$ws.Cells.Item(192, 7).Value = 360

$ws.Cells.Item(193, 1).Value = "Year 5"
$ws.Cells.Item(193, 2).Value = "B2-C1"
$ws.Cells.Item(193, 3).Value = "neurology"
$ws.Cells.Item(193, 4).Value = "'7"
$ws.Cells.Item(193, 5).Value = "'29/12/2025"
$ws.Cells.Item(193, 6).Value = "'09:00:00"
$ws.Cells.Item(193, 7).Value = 360

$ws.Cells.Item(194, 1).Value = "Year 5"
$ws.Cells.Item(194, 2).Value = "B2-C1"
$ws.Cells.Item(194, 3).Value = "neurology"
$ws.Cells.Item(194, 4).Value = "'8"
$ws.Cells.Item(194, 5).Value = "'30/12/2025"
$ws.Cells.Item(194, 6).Value = "'09:00:00"
$ws.Cells.Item(194, 7).Value = 360

$ws.Cells.Item(195, 1).Value = "Year 5"
$ws.Cells.Item(195, 2).Value = "B2-C1"
$ws.Cells.Item(195, 3).Value = "physical medicine"
$ws.Cells.Item(195, 4).Value = "'1"
$ws.Cells.Item(195, 5).Value = "'24/12/2025"
$ws.Cells.Item(195, 6).Value = "'09:00:00"
$ws.Cells.Item(195, 7).Value = 360

$ws.Cells.Item(196, 1).Value = "Year 5"
$ws.Cells.Item(196, 2).Value = "B2-C1"
$ws.Cells.Item(196, 3).Value = "physical medicine"
$ws.Cells.Item(196, 4).Value = "'2"
$ws.Cells.Item(196, 5).Value = "'31/12/2025"
$ws.Cells.Item(196, 6).Value = "'09:00:00"
$ws.Cells.Item(196, 7).Value = 360

$ws.Cells.Item(197, 1).Value = "Year 5"
$ws.Cells.Item(197, 2).Value = "B2-C1"
$ws.Cells.Item(197, 3).Value = "rheumatology"
$ws.Cells.Item(197, 4).Value = "'1"
$ws.Cells.Item(197, 5).Value = "'06/12/2025"
$ws.Cells.Item(197, 6).Value = "'09:00:00"
$ws.Cells.Item(197, 7).Value = 360

$ws.Cells.Item(198, 1).Value = "Year 5"
$ws.Cells.Item(198, 2).Value = "B2-C1"
$ws.Cells.Item(198, 3).Value = "rheumatology"
$ws.Cells.Item(198, 4).Value = "'2"
$ws.Cells.Item(198, 5).Value = "'07/12/2025"
$ws.Cells.Item(198, 6).Value = "'09:00:00"
$ws.Cells.Item(198, 7).Value = 360

$ws.Cells.Item(199, 1).Value = "Year 5"
$ws.Cells.Item(199, 2).Value = "B2-C1"
$ws.Cells.Item(199, 3).Value = "rheumatology"
$ws.Cells.Item(199, 4).Value = "'3"
$ws.Cells.Item(199, 5).Value = "'08/12/2025"
$ws.Cells.Item(199, 6).Value = "'09:00:00"
$ws.Cells.Item(199, 7).Value = 360

$ws.Cells.Item(200, 1).Value = "Year 5"
$ws.Cells.Item(200, 2).Value = "B2-C1"
$ws.Cells.Item(200, 3).Value = "rheumatology"
$ws.Cells.Item(200, 4).Value = "'4"
$ws.Cells.Item(200, 5).Value = "'09/12/2025"
$ws.Cells.Item(200, 6).Value = "'09:00:00"
$ws.Cells.Item(200, 7).Value = 360

$ws.Cells.Item(201, 1).Value = "Year 5"
$ws.Cells.Item(201, 2).Value = "B2-C1"
$ws.Cells.Item(201, 3).Value = "rheumatology"
$ws.Cells.Item(201, 4).Value = "'5"
$ws.Cells.Item(201, 5).Value = "'10/12/2025"
$ws.Cells.Item(201, 6).Value = "'09:00:00"
$ws.Cells.Item(201, 7).Value = 360

# Step 2: copy number-format/fill/font/alignment (but not values) from the last two
# existing rows (160-161), which have the same alternating style pair used throughout
# the table, onto the newly added rows 162-201. PasteSpecial(xlPasteFormats) only touches
# formatting, so the literal text values set above survive untouched.
$fmtSrc = $ws.Range("A160:G161")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A162:G201")
$fmtDst.PasteSpecial(-4122)

$excel.CutCopyMode = 0